$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H61").Value = 1772.4
$ws1.Range("I61").Value = 1772.4
$ws1.Range("J61").Value = 0
$ws1.Range("K61").Value = 5317.200000000001
$ws1.Range("L61").Value = 0
$ws1.Range("M61").Value = -5145.200000000001
$ws1.Range("N61").ClearContents()
$ws1.Range("H74").Value = 8022.84
$ws1.Range("I74").Value = 6228.846
$ws1.Range("K74").Value = 6228.846
$ws1.Range("M74").Value = -5292.846
$ws1.Range("H77").Value = 8022.84
$ws1.Range("I77").Value = 6228.846
$ws1.Range("K77").Value = 31144.23
$ws1.Range("M77").Value = -26464.23
$ws1.Range("H80").Value = 295.42105
$ws1.Range("I80").Value = 302.42856
$ws1.Range("J80").Value = 291.33334
$ws1.Range("K80").Value = 907.28568
$ws1.Range("L80").Value = 874.0000200000001
$ws1.Range("M80").Value = 90.71432000000004
$ws1.Range("N80").Value = -2870.00002
$ws1.Range("H83").Value = 295.42105
$ws1.Range("I83").Value = 302.42856
$ws1.Range("J83").Value = 291.33334
$ws1.Range("K83").Value = 2721.85704
$ws1.Range("L83").Value = 2622.00006
$ws1.Range("M83").Value = 2270.14296
$ws1.Range("N83").Value = -12606.00006
$ws1.Range("H87").Value = 80487.25
$ws1.Range("J87").Value = 90649.664
$ws1.Range("L87").Value = 90649.664
$ws1.Range("N87").Value = -93145.664
$ws1.Range("H90").Value = 80487.25
$ws1.Range("J90").Value = 90649.664
$ws1.Range("L90").Value = 271948.992
$ws1.Range("N90").Value = -284428.992
$ws1.Range("H125").Value = 2846.2
$ws1.Range("I125").Value = 2391.9
$ws1.Range("J125").Value = 3300.5
$ws1.Range("K125").Value = 21527.1
$ws1.Range("L125").Value = 29704.5
$ws1.Range("M125").Value = -19067.1
$ws1.Range("N125").Value = -34624.5
$ws1.Range("H132").Value = 384040.7
$ws1.Range("I132").Value = 419265.97
$ws1.Range("K132").Value = 1257797.91
$ws1.Range("M132").Value = -1255267.91
$ws1.Range("H138").Value = 3209.322
$ws1.Range("I138").Value = 2494.4285
$ws1.Range("K138").Value = 7483.2855
$ws1.Range("M138").Value = -2343.2855

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H45").Value = 3282.8096
$ws2.Range("I45").Value = 3049.6843
$ws2.Range("J45").Value = 5497.5
$ws2.Range("K45").Value = 3049.6843
$ws2.Range("L45").Value = 5497.5
$ws2.Range("M45").Value = -2672.6843
$ws2.Range("N45").Value = -6251.5
$ws2.Range("H61").Value = 5760.3335
$ws2.Range("J61").Value = 12308.363
$ws2.Range("L61").Value = 12308.363
$ws2.Range("N61").Value = -12732.363
$ws2.Range("H74").Value = 2608811
$ws2.Range("I74").Value = 3575260
$ws2.Range("K74").Value = 3575260
$ws2.Range("M74").Value = -3574386
$ws2.Range("H77").Value = 2608811
$ws2.Range("I77").Value = 3575260
$ws2.Range("K77").Value = 17876300
$ws2.Range("M77").Value = -17871932
$ws2.Range("H132").Value = 936543.9399999999
$ws2.Range("I132").Value = 1458080.2
$ws2.Range("K132").Value = 4374240.6
$ws2.Range("M132").Value = -4371710.6
$ws2.Range("H136").Value = 5760.3335
$ws2.Range("J136").Value = 12308.363
$ws2.Range("L136").Value = 36925.089
$ws2.Range("N136").Value = -42025.089

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H94").Value = 824.76666
$ws3.Range("I94").Value = 680.64703
$ws3.Range("K94").Value = 680.64703
$ws3.Range("M94").Value = -229.64703
$ws3.Range("H134").Value = 554948.3
$ws3.Range("I134").Value = 655757.1
$ws3.Range("J134").Value = 7700.2856
$ws3.Range("K134").Value = 1967271.3
$ws3.Range("L134").Value = 23100.8568
$ws3.Range("M134").Value = -1964736.3
$ws3.Range("N134").Value = -28170.8568

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H31").Value = 10106.235
$ws4.Range("I31").Value = 2128.8333
$ws4.Range("K31").Value = 2128.8333
$ws4.Range("M31").Value = -1833.8333
$ws4.Range("H34").Value = 10106.235
$ws4.Range("I34").Value = 2128.8333
$ws4.Range("K34").Value = 2128.8333
$ws4.Range("M34").Value = -1926.8333
$ws4.Range("H58").Value = 688254.25
$ws4.Range("I58").Value = 824371.9
$ws4.Range("J58").Value = 7666.3335
$ws4.Range("K58").Value = 824371.9
$ws4.Range("L58").Value = 7666.3335
$ws4.Range("M58").Value = -824168.9
$ws4.Range("N58").Value = -8072.3335
$ws4.Range("H70").Value = 75525
$ws4.Range("J70").Value = 75394
$ws4.Range("L70").Value = 75394
$ws4.Range("N70").Value = -76024
$ws4.Range("H73").Value = 75525
$ws4.Range("J73").Value = 75394
$ws4.Range("L73").Value = 75394
$ws4.Range("N73").Value = -77578
$ws4.Range("H132").Value = 6261372.5
$ws4.Range("I132").Value = 12530.156
$ws4.Range("J132").Value = 31256742
$ws4.Range("K132").Value = 37590.468
$ws4.Range("L132").Value = 93770226
$ws4.Range("M132").Value = -35060.468
$ws4.Range("N132").Value = -93775286
$ws4.Range("H134").Value = 13885
$ws4.Range("I134").Value = 12259.735
$ws4.Range("K134").Value = 36779.205
$ws4.Range("M134").Value = -34244.205
$ws4.Range("H136").Value = 688254.25
$ws4.Range("I136").Value = 824371.9
$ws4.Range("J136").Value = 7666.3335
$ws4.Range("K136").Value = 2473115.7
$ws4.Range("L136").Value = 22999.0005
$ws4.Range("M136").Value = -2470565.7
$ws4.Range("N136").Value = -28099.0005
$ws4.Range("H141").Value = 424634
$ws4.Range("J141").Value = 550178.8
$ws4.Range("L141").Value = 550178.8
$ws4.Range("N141").Value = -560538.8

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H25").Value = 0
$ws5.Range("J25").Value = 0
$ws5.Range("L25").Value = 0
$ws5.Range("N25").ClearContents()
$ws5.Range("H30").Value = 0
$ws5.Range("J30").Value = 0
$ws5.Range("L30").Value = 0
$ws5.Range("N30").ClearContents()
$ws5.Range("H39").Value = 425
$ws5.Range("J39").Value = 0
$ws5.Range("L39").Value = 0
$ws5.Range("N39").ClearContents()
$ws5.Range("H41").Value = 149
$ws5.Range("I41").Value = 98
$ws5.Range("J41").Value = 200
$ws5.Range("K41").Value = 294
$ws5.Range("L41").Value = 600
$ws5.Range("M41").Value = 44
$ws5.Range("N41").Value = -1276
$ws5.Range("H55").Value = 4901.3335
$ws5.Range("J55").Value = 9900
$ws5.Range("L55").Value = 29700
$ws5.Range("N55").Value = -30054
$ws5.Range("H130").Value = 2200
$ws5.Range("I130").Value = 2200
$ws5.Range("K130").Value = 6600
$ws5.Range("M130").Value = -1580
$ws5.Range("H131").Value = 9506.759
$ws5.Range("I131").Value = 912.6667
$ws5.Range("J131").Value = 11748.695
$ws5.Range("K131").Value = 2738.0001
$ws5.Range("L131").Value = 35246.085
$ws5.Range("M131").Value = 2301.9999
$ws5.Range("N131").Value = -45326.085

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H80").Value = 151735.67
$ws6.Range("I80").Value = 230814.77
$ws6.Range("J80").Value = 6757.3335
$ws6.Range("K80").Value = 230814.77
$ws6.Range("L80").Value = 6757.3335
$ws6.Range("M80").Value = -229816.77
$ws6.Range("N80").Value = -8753.333500000001
$ws6.Range("H83").Value = 151735.67
$ws6.Range("I83").Value = 230814.77
$ws6.Range("J83").Value = 6757.3335
$ws6.Range("K83").Value = 1154073.85
$ws6.Range("L83").Value = 33786.6675
$ws6.Range("M83").Value = -1149081.85
$ws6.Range("N83").Value = -43770.6675
$ws6.Range("H97").Value = 1314.8966
$ws6.Range("I97").Value = 631.6
$ws6.Range("J97").Value = 2047
$ws6.Range("K97").Value = 631.6
$ws6.Range("L97").Value = 2047
$ws6.Range("M97").Value = -135.6
$ws6.Range("N97").Value = -3039
$ws6.Range("H113").Value = 1682.6666
$ws6.Range("I113").Value = 1659.4
$ws6.Range("J113").Value = 1799
$ws6.Range("K113").Value = 1659.4
$ws6.Range("L113").Value = 1799
$ws6.Range("M113").Value = 510.5999999999999
$ws6.Range("N113").Value = -6139
$ws6.Range("H132").Value = 5212.846
$ws6.Range("I132").Value = 3217.625
$ws6.Range("K132").Value = 9652.875
$ws6.Range("M132").Value = -7122.875

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H132").Value = 1388718.2
$ws7.Range("I132").Value = 1824366.1
$ws7.Range("J132").Value = 9166.5
$ws7.Range("K132").Value = 5473098.300000001
$ws7.Range("L132").Value = 27499.5
$ws7.Range("M132").Value = -5470568.300000001
$ws7.Range("N132").Value = -32559.5

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H81").Value = 1790.4
$ws8.Range("I81").Value = 1987.5
$ws8.Range("J81").Value = 1002
$ws8.Range("K81").Value = 3975
$ws8.Range("L81").Value = 2004
$ws8.Range("M81").Value = -2914
$ws8.Range("N81").Value = -4126
$ws8.Range("H84").Value = 1790.4
$ws8.Range("I84").Value = 1987.5
$ws8.Range("J84").Value = 1002
$ws8.Range("K84").Value = 19875
$ws8.Range("L84").Value = 10020
$ws8.Range("M84").Value = -14571
$ws8.Range("N84").Value = -20628
$ws8.Range("H122").Value = 2259.7222
$ws8.Range("I122").Value = 1953.7931
$ws8.Range("J122").Value = 3527.1428
$ws8.Range("K122").Value = 5861.379300000001
$ws8.Range("L122").Value = 10581.4284
$ws8.Range("M122").Value = -3411.379300000001
$ws8.Range("N122").Value = -15481.4284
$ws8.Range("H132").Value = 28459426
$ws8.Range("I132").Value = 2882704
$ws8.Range("K132").Value = 8648112
$ws8.Range("M132").Value = -8645582
